# Refine format of code, comments and excel templates
# The header cell B1 on the "TEM" sheet previously held the placeholder
# label "TEMPLATE寄存器" (merged across B1:G1). Update it to the shorter
# "TEM" label to match the renamed template/sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "TEM"
